# "Generate Report for Handback"
#
# This localization-status report previously showed every file as
# "Ready for handoff" with placeholder handback timestamps
# (0001-01-01 00:00:00). This script updates the report to reflect that
# the zh-cn and de-de handbacks have happened: the status text changes,
# each language table gets its own "Latest Handback DateTime", and two
# new columns (Latest Target File / Latest Handback File) are populated
# with hyperlinked file names for every row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shared by the Overview summary columns (B/C) and the
#    per-language Status column (C) on both language sheets.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Latest Handback DateTime (column H): replace the 0001-01-01
#    placeholder with real handback timestamps - distinct per language.
# ---------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-11 20:44:46"
$wsZhCn.Range("H3").Value = "2016-03-11 20:44:46"

$wsDeDe.Range("H2").Value = "2016-03-11 20:44:52"
$wsDeDe.Range("H3").Value = "2016-03-11 20:44:52"

# ---------------------------------------------------------------------
# 3. Populate "Latest Target File" (F) and "Latest Handback File" (G)
#    with hyperlinked file names, for both rows of both language sheets.
# ---------------------------------------------------------------------

# zh-cn
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e19a416bef887603d1b580c59daca17e18242bec/e2e/59d3de95-70ff-4265-b60a-42f586f1ad0c.md",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.md") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/351e038cdc52d48117044d7443af796a00b0253e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e19a416bef887603d1b580c59daca17e18242bec/e2e/59d3de95-70ff-4265-b60a-42f586f1ad0c.md",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.md") | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/351e038cdc52d48117044d7443af796a00b0253e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.zh-cn.xlf") | Out-Null

# de-de
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e19a416bef887603d1b580c59daca17e18242bec/e2e/59d3de95-70ff-4265-b60a-42f586f1ad0c.md",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.md") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/961473f475ed9a076f76e1d60e2b772e04e38a01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e19a416bef887603d1b580c59daca17e18242bec/e2e/59d3de95-70ff-4265-b60a-42f586f1ad0c.md",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.md") | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/961473f475ed9a076f76e1d60e2b772e04e38a01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf",
    "",
    "",
    "59d3de95-70ff-4265-b60a-42f586f1ad0c.6184a512d1b8dac99c27845d59ea3208e07ff84c.de-de.xlf") | Out-Null

Write-Output "Handback report generated."
